# Applies the "Bolão independência - 30.xlsx" revision:
#   - sheet renamed from "10 cotas(3 Jogos)" to "12 cotas"
#   - header cell A1 re-affirmed with its bold/centered formatting
#   - a new (empty, underline-formatted) cell C33 added below the table,
#     which extends the sheet's used range to A1:O33

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab: "10 cotas(3 Jogos)" -> "12 cotas"
$ws.Name = "12 cotas"

# Keep the "Bola 1" header (A1) bold/centered, matching the rest of row 1.
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Name = "Calibri"
$ws.Range("A1").Font.Size = 11
$ws.Range("A1").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A1").VerticalAlignment = -4160     # xlTop

# Add the new trailing cell C33 (row 32 stays blank), formatted with an
# underline - this grows the sheet's dimension to A1:O33.
$ws.Range("C33").Font.Underline = $true
